$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J8").Value = 4.5
$ws.Range("J9").Value = 5
$ws.Range("G10").Value = 5.5
$ws.Range("J11").Value = 4.5
$ws.Range("G12").Value = 5.5
$ws.Range("G13").Value = 5.5
$ws.Range("J13").Value = $null
$ws.Range("G14").Value = 5.5
$ws.Range("J15").Value = 5

$ws.Range("M19").Value = (Get-Date -Year 2022 -Month 11 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M21").Value = "Al Asfar Ahmad"
$ws.Range("M23").Value = "XCL"

$ws.Range("M23:O23").Select() | Out-Null
